$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove rows 5 and 6 (Happy Helper and Nustay entries dropped from this extract)
$ws.Range("A5:A6").EntireRow.Delete()

# Clear cells that no longer have values after the update
$ws.Range("T2").ClearContents()
$ws.Range("D4").ClearContents()
$ws.Range("T4").ClearContents()

# Update row 2 (company id/financials refreshed)
$ws.Range("A2").Value = 'Denmark'
$ws.Range("B2").Value = "'2"
$ws.Range("C2").Value = 'Retail (Online)'
$ws.Range("D2").Value = -0.012
$ws.Range("G2").Value = -0.01147982062780269
$ws.Range("H2").Value = -0.08370702541106129
$ws.Range("I2").Value = -0.129745889387145
$ws.Range("J2").Value = -0.129745889387145
$ws.Range("K2").Value = -5.464
$ws.Range("L2").Value = -0.1633482810164425
$ws.Range("M2").Value = 0
$ws.Range("N2").Value = 0
$ws.Range("O2").Value = -0
$ws.Range("P2").Value = 0
$ws.Range("Q2").Value = 0
$ws.Range("R2").Value = -0
$ws.Range("S2").Value = 0
$ws.Range("U2").Value = 4.39
$ws.Range("V2").Value = 0.3395204949729311
$ws.Range("W2").Value = -0.9956823318525445
$ws.Range("X2").Value = 0.2238688589878753
$ws.Range("Y2").Value = -1.21955119084042
$ws.Range("Z2").Value = 2.445711778898881
$ws.Range("AA2").Value = -0.4163962310251287
$ws.Range("AB2").Value = 0.06990772922347488
$ws.Range("AC2").Value = -0.4863039602486036
$ws.Range("AD2").Value = 30.639
$ws.Range("AE2").Value = 0
$ws.Range("AF2").Value = 30.639
$ws.Range("AG2").Value = 26.249
$ws.Range("AH2").Value = 0.7032293603250017
$ws.Range("AI2").Value = 0.9773517496570864
$ws.Range("AJ2").Value = 0.6699762627938437
$ws.Range("AK2").Value = 0.9736637115620015
$ws.Range("AL2").Value = 1.819
$ws.Range("AM2").Value = 1.485
$ws.Range("AN2").Value = -5.264432989690722
$ws.Range("AO2").Value = -2.385926333150083
$ws.Range("AP2").Value = -4.510137457044674
$ws.Range("AQ2").Value = -2.922558922558923

# Update row 3 (now Lauritz.com Group A/S)
$ws.Range("A3").Value = 'Denmark'
$ws.Range("B3").Value = 'Lauritz.com Group A/S (OM:LAUR)'
$ws.Range("C3").Value = 'Retail (Online)'
$ws.Range("D3").Value = -0.012
$ws.Range("G3").Value = 0.0003095975232198145
$ws.Range("H3").Value = -0.05572755417956657
$ws.Range("I3").Value = -0.1018575851393189
$ws.Range("J3").Value = -0.1018575851393189
$ws.Range("K3").Value = -4.6
$ws.Range("L3").Value = -0.1424148606811146
$ws.Range("M3").Value = -0
$ws.Range("N3").Value = -0
$ws.Range("O3").Value = 0
$ws.Range("P3").Value = -0
$ws.Range("Q3").Value = -0
$ws.Range("R3").Value = 0
$ws.Range("S3").Value = 0
$ws.Range("U3").Value = 3.54
$ws.Range("V3").Value = 0.735966735966736
$ws.Range("W3").Value = -1.684981684981685
$ws.Range("X3").Value = 0.3733337216043537
$ws.Range("Y3").Value = -2.058315406586039
$ws.Range("Z3").Value = 2.741935483870967
$ws.Range("AA3").Value = -0.2792869269949066
$ws.Range("AB3").Value = 0.06879024923852597
$ws.Range("AC3").Value = -0.3480771762334326
$ws.Range("AD3").Value = 30.1
$ws.Range("AE3").Value = 0
$ws.Range("AF3").Value = 30.1
$ws.Range("AG3").Value = 26.56
$ws.Range("AH3").Value = 0.8622171297622457
$ws.Range("AI3").Value = 1.052079692415239
$ws.Range("AJ3").Value = 0.846668791839337
$ws.Range("AK3").Value = 1.059433585959314
$ws.Range("AL3").Value = 1.75
$ws.Range("AM3").Value = 1.416
$ws.Range("AN3").Value = -6.244813278008299
$ws.Range("AO3").Value = -1.88
$ws.Range("AP3").Value = -5.510373443983402
$ws.Range("AQ3").Value = -2.323446327683616

# Update row 4 (now Happy Helper A/S)
$ws.Range("A4").Value = 'Denmark'
$ws.Range("B4").Value = 'Happy Helper A/S (CPSE:HAPPY)'
$ws.Range("C4").Value = 'Retail (Online)'
$ws.Range("G4").Value = -0.3426086956521739
$ws.Range("H4").Value = -0.8695652173913044
$ws.Range("I4").Value = -0.9130434782608696
$ws.Range("J4").Value = -0.9130434782608696
$ws.Range("K4").Value = -0.864
$ws.Range("L4").Value = -0.751304347826087
$ws.Range("M4").Value = -0
$ws.Range("N4").Value = -0
$ws.Range("O4").Value = 0
$ws.Range("P4").Value = -0
$ws.Range("Q4").Value = -0
$ws.Range("R4").Value = 0
$ws.Range("S4").Value = 0
$ws.Range("U4").Value = 0.85
$ws.Range("V4").Value = 0.104679802955665
$ws.Range("W4").Value = -0.3063829787234043
$ws.Range("X4").Value = 0.07440399637139676
$ws.Range("Y4").Value = -0.380786975094801
$ws.Range("Z4").Value = 0.606220347917765
$ws.Range("AA4").Value = -0.5535055350553507
$ws.Range("AB4").Value = 0.0710252092084238
$ws.Range("AC4").Value = -0.6245307442637745
$ws.Range("AD4").Value = 0.539
$ws.Range("AE4").Value = 0
$ws.Range("AF4").Value = 0.539
$ws.Range("AG4").Value = -0.3109999999999999
$ws.Range("AH4").Value = 0.06224737267582863
$ws.Range("AI4").Value = 0.1967871485943775
$ws.Range("AJ4").Value = -0.03982584197720578
$ws.Range("AK4").Value = -0.1646373742721016
$ws.Range("AL4").Value = 0.06900000000000001
$ws.Range("AM4").Value = 0.06900000000000001
$ws.Range("AN4").Value = -0.539
$ws.Range("AO4").Value = -15.21739130434783
$ws.Range("AP4").Value = 0.3109999999999999
$ws.Range("AQ4").Value = -15.21739130434783
